$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Y-Strainer"
$ws.Range("D3").Value = "Ever"
$ws.Range("B4").Value = "150 lbs Rating, Flange Type, Cast Iron, 8`""
